$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("Q3").Formula  = "'18190.34"
$ws.Range("Q4").Formula  = "'17030.91"
$ws.Range("Q5").Formula  = "'24652.88"
$ws.Range("Q6").Formula  = "'14160.39"
$ws.Range("Q7").Formula  = "'21043.53"
$ws.Range("Q8").Formula  = "'10547.53"
$ws.Range("Q9").Formula  = "'14411.72"
$ws.Range("Q10").Formula = "'22075.29"
$ws.Range("Q11").Formula = "'19241.00"
$ws.Range("Q12").Formula = "'17032.48"
$ws.Range("Q13").Formula = "'91878.19"
$ws.Range("Q14").Formula = "'15225.24"
